$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A usernames to the final de-duplicated / reordered list
$ws.Range("A1").Value = 'katyperry'
$ws.Range("A2").Value = 'justinbieber'
$ws.Range("A3").Value = 'BarackObama'
$ws.Range("A4").Value = 'taylorswift13'
$ws.Range("A5").Value = 'ladygaga'
$ws.Range("A6").Value = 'rihanna'
$ws.Range("A7").Value = 'jtimberlake'
$ws.Range("A8").Value = 'TheEllenShow'
$ws.Range("A9").Value = 'britneyspears'
$ws.Range("A10").Value = 'Cristiano'
$ws.Range("A11").Value = 'KimKardashian'
$ws.Range("A12").Value = 'JLO'
$ws.Range("A13").Value = 'shakira'
$ws.Range("A14").Value = 'selenagomez'
$ws.Range("A15").Value = 'ArianaGrande'
$ws.Range("A16").Value = 'ddlovato'
$ws.Range("A17").Value = 'Oprah'
$ws.Range("A18").Value = 'Pink'
$ws.Range("A19").Value = 'jimmyfallon'
$ws.Range("A20").Value = 'Harry_Styles'
$ws.Range("A21").Value = 'onedirection'
$ws.Range("A22").Value = 'Drake'
$ws.Range("A23").Value = 'LilTunechi'
$ws.Range("A24").Value = 'KAKA'
$ws.Range("A25").Value = 'BillGates'
$ws.Range("A26").Value = 'NiallOfficial'
$ws.Range("A27").Value = 'aliciakeys'
$ws.Range("A28").Value = 'KingJames'
$ws.Range("A29").Value = 'BrunoMars'
$ws.Range("A30").Value = 'pitbull'
$ws.Range("A31").Value = 'Real_Liam_Payne'
$ws.Range("A32").Value = 'wizkhalifa'
$ws.Range("A33").Value = 'MileyCyrus'
$ws.Range("A34").Value = 'Louis_Tomlinson'
$ws.Range("A35").Value = 'KevinHart4real'
$ws.Range("A36").Value = 'Eminem'
$ws.Range("A37").Value = 'NICKIMINAJ'
$ws.Range("A38").Value = 'AvrilLavigne'
$ws.Range("A39").Value = 'neymarjr'
$ws.Range("A40").Value = 'davidguetta'
$ws.Range("A41").Value = 'danieltosh'
$ws.Range("A42").Value = 'aplusk'
$ws.Range("A43").Value = 'ConanOBrien'
$ws.Range("A44").Value = 'ActuallyNPH'
$ws.Range("A45").Value = 'MariahCarey'
$ws.Range("A46").Value = 'SrBachchan'
$ws.Range("A47").Value = 'coldplay'
$ws.Range("A48").Value = 'xtina'
$ws.Range("A49").Value = 'koutneykardash'
$ws.Range("A50").Value = 'zaynmalik'
$ws.Range("A51").Value = 'JimCarrey'
$ws.Range("A52").Value = 'khloekardashian'
$ws.Range("A53").Value = 'chrisbrown'
$ws.Range("A54").Value = 'Beyonce'
$ws.Range("A55").Value = 'edsheeran'
$ws.Range("A56").Value = 'RyanSeacrest'
$ws.Range("A57").Value = 'iamsrk'
$ws.Range("A58").Value = 'ParisHilton'
$ws.Range("A59").Value = 'agnezmo'
$ws.Range("A60").Value = 'iamwill'
$ws.Range("A61").Value = 'ivetesangalo'
$ws.Range("A62").Value = 'aamir_khan'
$ws.Range("A63").Value = 'LeoDiCaprio'
$ws.Range("A64").Value = 'narendramodi'
$ws.Range("A65").Value = 'ashleytisdale'
$ws.Range("A66").Value = 'kanyewest'
$ws.Range("A67").Value = 'tyrabanks'
$ws.Range("A68").Value = 'AlejandroSanz'
$ws.Range("A69").Value = 'blakeshelton'
$ws.Range("A70").Value = 'BeingSalmanKhan'
$ws.Range("A71").Value = 'SnoopDogg'
$ws.Range("A72").Value = '10Ronaldinho'
$ws.Range("A73").Value = 'ricky_martin'
$ws.Range("A74").Value = 'SimonCowell'
$ws.Range("A75").Value = 'MohamadAlarefe'
$ws.Range("A76").Value = 'charliesheen'
$ws.Range("A77").Value = 'ClaudiaLeitte'
$ws.Range("A78").Value = 'DalaiLama'
$ws.Range("A79").Value = 'maroon5'
$ws.Range("A80").Value = 'KendallJenner'
$ws.Range("A81").Value = 'andresiniesta8'
$ws.Range("A82").Value = 'carlyraejepsen'
$ws.Range("A83").Value = 'ZacEfron'
$ws.Range("A84").Value = 'marcosmion'
$ws.Range("A85").Value = 'LucianoHuck'
$ws.Range("A86").Value = '3gerardpique'
$ws.Range("A87").Value = 'radityadika'
$ws.Range("A88").Value = 'deepikapadukone'
$ws.Range("A89").Value = 'juanes'
$ws.Range("A90").Value = 'iamdiddy'
$ws.Range("A91").Value = 'KDtrey5'
$ws.Range("A92").Value = 'Ludacris'
$ws.Range("A93").Value = 'paurubio'
$ws.Range("A94").Value = 'DaniloGentilo'
$ws.Range("A95").Value = 'kelly_Clarkson'
$ws.Range("A96").Value = 'Usher'
$ws.Range("A97").Value = 'paulocoelho'

# Rows 98-100 lost their username after removing 3 duplicate entries from the list
$ws.Range("A98:A100").ClearContents()

# Column A best-fit width for the (now shorter) usernames
$ws.Columns("A").ColumnWidth = 14

# Restore the view: zoomed in, scrolled down near the bottom of the list,
# with A93:A97 selected (the re-ordered/consent-changed rows)
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 125
$win.ScrollRow = 88
$win.ScrollColumn = 1
$ws.Range("A93:A97").Select()
